$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.114708065986633
$ws.Range("B1").Value = 2.809435844421387
$ws.Range("C1").Value = 8.789849281311035
$ws.Range("D1").Value = 2.031454801559448
$ws.Range("E1").Value = 1.137068510055542
